$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100
$ws.Range("B100").Value2 = 6867460
$ws.Range("F100").Value = "FC Vion Zlate Moravce"
$ws.Range("G100").Value = "FC Kosice"
$ws.Range("H100").Value2 = 1
$ws.Range("I100").Value2 = 1
$ws.Range("J100").Value = "D"
$ws.Range("K100").Value2 = 2.3
$ws.Range("L100").Value2 = 3.3
$ws.Range("M100").Value2 = 2.875
$ws.Range("N100").Value2 = 2.75
$ws.Range("O100").Value2 = 3.1
$ws.Range("P100").Value2 = 2.75
$ws.Range("Q100").Value2 = 0
$ws.Range("R100").Value2 = 1.875
$ws.Range("S100").Value2 = 1.925
$ws.Range("T100").Value2 = 2
$ws.Range("U100").Value2 = 1.85
$ws.Range("V100").Value2 = 1.95
$ws.Range("W100").Value2 = -1
$ws.Range("X100").Value2 = 2.1
$ws.Range("Y100").Value2 = -1
$ws.Range("Z100").Value2 = 0
$ws.Range("AA100").Value2 = -0
$ws.Range("AB100").Value2 = 0
$ws.Range("AC100").Value2 = -0

# Row 101
$ws.Range("B101").Value2 = 6867461
$ws.Range("F101").Value = "FK Zeleziarne Podbrezova"
$ws.Range("G101").Value = "Slovan Bratislava"
$ws.Range("H101").Value2 = 0
$ws.Range("I101").Value2 = 6
$ws.Range("J101").Value = "A"
$ws.Range("K101").Value2 = 3.25
$ws.Range("L101").Value2 = 3.5
$ws.Range("M101").Value2 = 2
$ws.Range("N101").Value2 = 3.5
$ws.Range("O101").Value2 = 3.75
$ws.Range("P101").Value2 = 2
$ws.Range("Q101").Value2 = 0.5
$ws.Range("R101").Value2 = 1.8
$ws.Range("S101").Value2 = 2
$ws.Range("T101").Value2 = 3
$ws.Range("U101").Value2 = 1.95
$ws.Range("V101").Value2 = 1.85
$ws.Range("W101").Value2 = -1
$ws.Range("X101").Value2 = -1
$ws.Range("Y101").Value2 = 1
$ws.Range("Z101").Value2 = -1
$ws.Range("AA101").Value2 = 1
$ws.Range("AB101").Value2 = 0.95
$ws.Range("AC101").Value2 = -1

# Row 124
$ws.Range("B124").Value2 = 6867488
$ws.Range("F124").Value = "MSK Zilina"
$ws.Range("G124").Value = "MFK Zemplin Michalovce"
$ws.Range("H124").Value2 = 1
$ws.Range("I124").Value2 = 1
$ws.Range("J124").Value = "D"
$ws.Range("K124").Value2 = 1.45
$ws.Range("L124").Value2 = 4.333
$ws.Range("M124").Value2 = 5.75
$ws.Range("N124").Value2 = 1.45
$ws.Range("O124").Value2 = 4.5
$ws.Range("P124").Value2 = 6.5
$ws.Range("Q124").Value2 = -1.25
$ws.Range("R124").Value2 = 2
$ws.Range("S124").Value2 = 1.8
$ws.Range("T124").Value2 = 3
$ws.Range("U124").Value2 = 1.8
$ws.Range("V124").Value2 = 2
$ws.Range("W124").Value2 = -1
$ws.Range("X124").Value2 = 3.5
$ws.Range("Y124").Value2 = -1
$ws.Range("Z124").Value2 = -1
$ws.Range("AA124").Value2 = 0.8
$ws.Range("AB124").Value2 = -1
$ws.Range("AC124").Value2 = 1

# Row 125
$ws.Range("B125").Value2 = 6867489
$ws.Range("F125").Value = "FC Spartak Trnava"
$ws.Range("G125").Value = "Dukla Banska Bystrica"
$ws.Range("H125").Value2 = 2
$ws.Range("I125").Value2 = 0
$ws.Range("J125").Value = "H"
$ws.Range("K125").Value2 = 1.666
$ws.Range("L125").Value2 = 3.75
$ws.Range("M125").Value2 = 4.5
$ws.Range("N125").Value2 = 1.615
$ws.Range("O125").Value2 = 4
$ws.Range("P125").Value2 = 5.25
$ws.Range("Q125").Value2 = -0.75
$ws.Range("R125").Value2 = 1.775
$ws.Range("S125").Value2 = 2.025
$ws.Range("T125").Value2 = 2.5
$ws.Range("U125").Value2 = 1.825
$ws.Range("V125").Value2 = 1.975
$ws.Range("W125").Value2 = 0.615
$ws.Range("X125").Value2 = -1
$ws.Range("Y125").Value2 = -1
$ws.Range("Z125").Value2 = 0.7749999999999999
$ws.Range("AA125").Value2 = -1
$ws.Range("AB125").Value2 = -1
$ws.Range("AC125").Value2 = 0.9750000000000001

# Row 128
$ws.Range("B128").Value2 = 6867491
$ws.Range("F128").Value = "Dukla Banska Bystrica"
$ws.Range("G128").Value = "FC Kosice"
$ws.Range("H128").Value2 = 1
$ws.Range("I128").Value2 = 1
$ws.Range("J128").Value = "D"
$ws.Range("K128").Value2 = 1.571
$ws.Range("L128").Value2 = 4.2
$ws.Range("M128").Value2 = 5.25
$ws.Range("N128").Value2 = 1.571
$ws.Range("O128").Value2 = 4.2
$ws.Range("P128").Value2 = 5.5
$ws.Range("Q128").Value2 = -1
$ws.Range("R128").Value2 = 1.9
$ws.Range("S128").Value2 = 1.9
$ws.Range("T128").Value2 = 3
$ws.Range("U128").Value2 = 2
$ws.Range("V128").Value2 = 1.8
$ws.Range("W128").Value2 = -1
$ws.Range("X128").Value2 = 3.2
$ws.Range("Y128").Value2 = -1
$ws.Range("Z128").Value2 = -1
$ws.Range("AA128").Value2 = 0.8999999999999999
$ws.Range("AB128").Value2 = -1
$ws.Range("AC128").Value2 = 0.8

# Row 131
$ws.Range("B131").Value2 = 6867493
$ws.Range("F131").Value = "MFK Ruzomberok"
$ws.Range("G131").Value = "FK Zeleziarne Podbrezova"
$ws.Range("H131").Value2 = 2
$ws.Range("I131").Value2 = 1
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value2 = 3.3
$ws.Range("L131").Value2 = 3.3
$ws.Range("M131").Value2 = 2.2
$ws.Range("N131").Value2 = 2.9
$ws.Range("O131").Value2 = 3.4
$ws.Range("P131").Value2 = 2.4
$ws.Range("Q131").Value2 = 0.25
$ws.Range("R131").Value2 = 1.75
$ws.Range("S131").Value2 = 2.05
$ws.Range("T131").Value2 = 2.5
$ws.Range("U131").Value2 = 1.9
$ws.Range("V131").Value2 = 1.9
$ws.Range("W131").Value2 = 1.9
$ws.Range("X131").Value2 = -1
$ws.Range("Y131").Value2 = -1
$ws.Range("Z131").Value2 = 0.75
$ws.Range("AA131").Value2 = -1
$ws.Range("AB131").Value2 = 0.8999999999999999
$ws.Range("AC131").Value2 = -1

# Row 148
$ws.Range("B148").Value2 = 7911450
$ws.Range("F148").Value = "MFK Skalica"
$ws.Range("G148").Value = "MFK Zemplin Michalovce"
$ws.Range("H148").Value2 = 0
$ws.Range("I148").Value2 = 0
$ws.Range("J148").Value = "D"
$ws.Range("K148").Value2 = 2.3
$ws.Range("L148").Value2 = 3.25
$ws.Range("M148").Value2 = 3.1
$ws.Range("N148").Value2 = 2.2
$ws.Range("O148").Value2 = 3.3
$ws.Range("P148").Value2 = 3.3
$ws.Range("Q148").Value2 = -0.25
$ws.Range("R148").Value2 = 1.95
$ws.Range("S148").Value2 = 1.85
$ws.Range("T148").Value2 = 2.25
$ws.Range("U148").Value2 = 2
$ws.Range("V148").Value2 = 1.8
$ws.Range("W148").Value2 = -1
$ws.Range("X148").Value2 = 2.3
$ws.Range("Y148").Value2 = -1
$ws.Range("Z148").Value2 = -0.5
$ws.Range("AA148").Value2 = 0.425
$ws.Range("AB148").Value2 = -1
$ws.Range("AC148").Value2 = 0.8

# Row 149
$ws.Range("B149").Value2 = 7911478
$ws.Range("F149").Value = "FC Vion Zlate Moravce"
$ws.Range("G149").Value = "FC Kosice"
$ws.Range("H149").Value2 = 1
$ws.Range("I149").Value2 = 2
$ws.Range("J149").Value = "A"
$ws.Range("K149").Value2 = 2.5
$ws.Range("L149").Value2 = 3.2
$ws.Range("M149").Value2 = 2.8
$ws.Range("N149").Value2 = 2.6
$ws.Range("O149").Value2 = 3.1
$ws.Range("P149").Value2 = 2.875
$ws.Range("Q149").Value2 = 0
$ws.Range("R149").Value2 = 1.775
$ws.Range("S149").Value2 = 2.025
$ws.Range("T149").Value2 = 2.25
$ws.Range("U149").Value2 = 1.8
$ws.Range("V149").Value2 = 2
$ws.Range("W149").Value2 = -1
$ws.Range("X149").Value2 = -1
$ws.Range("Y149").Value2 = 1.875
$ws.Range("Z149").Value2 = -1
$ws.Range("AA149").Value2 = 1.025
$ws.Range("AB149").Value2 = 0.8
$ws.Range("AC149").Value2 = -1

# Row 152
$ws.Range("N152").Value2 = 4.5
$ws.Range("O152").Value2 = 4.2
$ws.Range("P152").Value2 = 1.615
$ws.Range("R152").Value2 = 1.775
$ws.Range("S152").Value2 = 2.025
$ws.Range("U152").Value2 = 1.875
$ws.Range("V152").Value2 = 1.925

# Row 153
$ws.Range("E153").Value2 = 45388.41666666666

# Row 154
$ws.Range("N154").Value2 = 2.9
$ws.Range("P154").Value2 = 2.3
$ws.Range("Q154").Value2 = 0.25
$ws.Range("R154").Value2 = 1.775
$ws.Range("S154").Value2 = 2.025

# Row 156
$ws.Range("N156").Value2 = 1.75
$ws.Range("O156").Value2 = 3.4
$ws.Range("P156").Value2 = 4.333
$ws.Range("Q156").Value2 = -0.75
$ws.Range("R156").Value2 = 2.025
$ws.Range("S156").Value2 = 1.775
$ws.Range("U156").Value2 = 1.8
$ws.Range("V156").Value2 = 2
